$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.628.94"
$ws.Range("D3").Value = "1.842.83"
$ws.Range("E3").Value = "  -1.32%  "
$ws.Range("E4").Value = "  -0.43%  "
$ws.Range("D5").Value = "'314.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("D6").Value = "'0.9996"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").Value = "'0.4253"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.89%  "
$ws.Range("D8").Value = "'0.3646"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'45.66"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("D10").Value = "'0.07277"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.14%  "
$ws.Range("D11").Value = "'0.8974"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.49%  "
$ws.Range("E12").Value = "  -3.78%  "
$ws.Range("D13").Value = "1.873.19"
$ws.Range("E13").Value = "  -3.98%  "
$ws.Range("D14").Value = "'5.389"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.04%  "
$ws.Range("D15").Value = "'6.559"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.46%  "
$ws.Range("D16").Value = "'0.06850"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").Value = "'1.000"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "'78.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.000008871"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.30%  "
$ws.Range("D20").Value = "'0.9994"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("D21").Value = "'15.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.35%  "
$ws.Range("D22").Value = "27.629.98"
$ws.Range("E22").Value = "  -2.26%  "
$ws.Range("D23").Value = "'4.975"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.08%  "
$ws.Range("D24").Value = "'10.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.78%  "
$ws.Range("B25").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D25").Value = "2.042.51"
$ws.Range("E25").Value = "  -5.22%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'2.044"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("D27").Value = "'154.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.48%  "
$ws.Range("D28").Value = "'18.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.83%  "
$ws.Range("D29").Value = "'5.247"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.27%  "
$ws.Range("D30").Value = "'1.837"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.29%  "
$ws.Range("D31").Value = "'111.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.58%  "
$ws.Range("D32").Value = "'0.08875"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("D33").Value = "'0.7769"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.57%  "
$ws.Range("D34").Value = "'4.571"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.42%  "
$ws.Range("D35").Value = "'2.947"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.34%  "
$ws.Range("D36").Value = "'1.098"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.17%  "
$ws.Range("D37").Value = "'0.9991"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.40%  "
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("E39").Value = "  -1.81%  "
$ws.Range("D40").Value = "'0.01927"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.35%  "
$ws.Range("D41").Value = "'2.763"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.54%  "
$ws.Range("D42").Value = "'0.5067"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.50%  "
$ws.Range("D43").Value = "'6.801"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.92%  "
$ws.Range("D44").Value = "'0.1640"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.00%  "
$ws.Range("D45").Value = "'8.231"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.46%  "
$ws.Range("D46").Value = "'0.06634"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.86%  "
$ws.Range("D47").Value = "'10.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.01%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'105.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.98%  "
$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").Value = "'0.4711"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.42%  "
$ws.Range("D50").Value = "'0.9992"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.25%  "
